$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.196.47'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '3.271.43'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.27'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.81'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("E8").Value = '  +0.55%  '
$ws.Range("E9").Value = '  -1.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.63'
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("E11").Value = '  -2.73%  '
$ws.Range("D12").Value = '3.840.68'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("E13").Value = '  +1.08%  '
$ws.Range("D14").Value = '68.168.55'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.36'
$ws.Range("E15").Value = '  -2.79%  '
$ws.Range("E16").Value = '  -1.70%  '
$ws.Range("D17").Value = '3.272.51'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("E18").Value = '  -2.17%  '
$ws.Range("E19").Value = '  -2.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '416.58'
$ws.Range("E20").Value = '  +5.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.53'
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.17'
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("E24").Value = '  -2.09%  '
$ws.Range("E25").Value = '  -1.81%  '
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("E27").Value = '  -4.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.95'
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.66'
$ws.Range("E30").Value = '  -1.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.43'
$ws.Range("E31").Value = '  -4.66%  '
$ws.Range("E32").Value = '  -2.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.84'
$ws.Range("E33").Value = '  -4.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '164.49'
$ws.Range("E34").Value = '  +1.36%  '
$ws.Range("E35").Value = '  -4.12%  '
$ws.Range("E36").Value = '  -3.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.62'
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.790'
$ws.Range("E38").Value = '  -3.77%  '
$ws.Range("E39").Value = '  -3.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.29'
$ws.Range("E40").Value = '  -3.80%  '
$ws.Range("D41").Value = '2.632.64'
$ws.Range("E41").Value = '  -0.81%  '
$ws.Range("E42").Value = '  -2.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.41'
$ws.Range("E43").Value = '  -3.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '335.05'
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.19'
$ws.Range("E45").Value = '  -4.81%  '
$ws.Range("E46").Value = '  -3.03%  '
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("E48").Value = '  -2.25%  '
$ws.Range("E49").Value = '  -1.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.73'
$ws.Range("E50").Value = '  -2.73%  '
$ws.Range("E51").Value = '  +0.05%  '
